# Weekly update: insert a new observation row for
# "Feria Lagunitas de Puerto Montt - Betarraga" ahead of the existing
# history (rows shift down by one; a new last row is created).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 176 (and everything below it) down by one row,
# inheriting formatting (incl. the date number-format on column D).
$ws.Rows("176:176").Insert()

# Populate the freshly inserted row 176 with this week's data.
$ws.Range("A176").Value = 4
$ws.Range("B176").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C176").Value = "Los Lagos"
$ws.Range("D176").Value = 44504
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 100114014
$ws.Range("G176").Value = "Betarraga"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 500
$ws.Range("K176").Value = 1000
$ws.Range("L176").Value = 1000
$ws.Range("M176").Value = 1000
$ws.Range("N176").Value = '$/paquete 5 unidades'
$ws.Range("O176").Value = "Región del Maule"
$ws.Range("P176").Value = 200
$ws.Range("Q176").Value = 5
$ws.Range("R176").Value = "Hortaliza"
